$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 253.38461
$ws.Range("I12").Value = 266.33334
$ws.Range("J12").Value = 98
$ws.Range("K12").Value = 266.33334
$ws.Range("L12").Value = 98
$ws.Range("M12").Value = -96.33334000000002
$ws.Range("N12").Value = -438
$ws.Range("H112").Value = 15627148
$ws.Range("J112").Value = 16131205
$ws.Range("L112").Value = 48393615
$ws.Range("N112").Value = -48395831
$ws.Range("H116").Value = 11443.444
$ws.Range("I116").Value = 10331.667
$ws.Range("K116").Value = 10331.667
$ws.Range("M116").Value = -6889.666999999999
$ws.Range("H123").Value = 58200
$ws.Range("J123").Value = 58200
$ws.Range("L123").Value = 58200
$ws.Range("N123").Value = -68000
$ws.Range("H125").Value = 3809.75
$ws.Range("I125").Value = 1500
$ws.Range("K125").Value = 13500
$ws.Range("M125").Value = -11040

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1314.4546
$ws.Range("I2").Value = 1181.4286
$ws.Range("K2").Value = 1181.4286
$ws.Range("M2").Value = -1068.4286
$ws.Range("H45").Value = 2175
$ws.Range("I45").Value = 1565.2222
$ws.Range("K45").Value = 1565.2222
$ws.Range("M45").Value = -1188.2222
$ws.Range("H88").Value = 50178
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 50178
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H102").Value = 171067.25
$ws.Range("I102").Value = 336166
$ws.Range("J102").Value = 5968.5
$ws.Range("K102").Value = 336166
$ws.Range("L102").Value = 5968.5
$ws.Range("M102").Value = -334544
$ws.Range("N102").Value = -9212.5
$ws.Range("H110").Value = 24822.334
$ws.Range("J110").Value = 2850
$ws.Range("L110").Value = 2850
$ws.Range("N110").Value = -6940
$ws.Range("H116").Value = 1314.4546
$ws.Range("I116").Value = 1181.4286
$ws.Range("K116").Value = 1181.4286
$ws.Range("M116").Value = 1112.5714
$ws.Range("H124").Value = 27057.25
$ws.Range("J124").Value = 27057.25
$ws.Range("L124").Value = 27057.25
$ws.Range("N124").Value = -36877.25
$ws.Range("H134").Value = 91063.10000000001
$ws.Range("J134").Value = 91181.22
$ws.Range("L134").Value = 91181.22
$ws.Range("N134").Value = -101321.22
$ws.Range("H140").Value = 113856.75
$ws.Range("J140").Value = 113856.75
$ws.Range("L140").Value = 113856.75
$ws.Range("N140").Value = -124216.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1314.4546
$ws.Range("I3").Value = 1181.4286
$ws.Range("K3").Value = 1181.4286
$ws.Range("M3").Value = -1067.4286
$ws.Range("H99").Value = 3496.2856
$ws.Range("I99").Value = 1817.8
$ws.Range("K99").Value = 1817.8
$ws.Range("M99").Value = -319.8
$ws.Range("H107").Value = 3400.6667
$ws.Range("I107").Value = 3216.2632
$ws.Range("K107").Value = 3216.2632
$ws.Range("M107").Value = -1296.2632
$ws.Range("H140").Value = 103713.57
$ws.Range("J140").Value = 103713.57
$ws.Range("L140").Value = 103713.57
$ws.Range("N140").Value = -114073.57

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4814.537
$ws.Range("I31").Value = 3072.3928
$ws.Range("J31").Value = 6690.6924
$ws.Range("K31").Value = 3072.3928
$ws.Range("L31").Value = 6690.6924
$ws.Range("M31").Value = -2777.3928
$ws.Range("N31").Value = -7280.6924
$ws.Range("H34").Value = 4814.537
$ws.Range("I34").Value = 3072.3928
$ws.Range("J34").Value = 6690.6924
$ws.Range("K34").Value = 3072.3928
$ws.Range("L34").Value = 6690.6924
$ws.Range("M34").Value = -2870.3928
$ws.Range("N34").Value = -7094.6924
$ws.Range("H58").Value = 5980081.5
$ws.Range("I58").Value = 1283.5807
$ws.Range("K58").Value = 1283.5807
$ws.Range("M58").Value = -1080.5807
$ws.Range("H86").Value = 11399.556
$ws.Range("J86").Value = 14999.667
$ws.Range("L86").Value = 14999.667
$ws.Range("N86").Value = -17245.667
$ws.Range("H89").Value = 11399.556
$ws.Range("J89").Value = 14999.667
$ws.Range("L89").Value = 74998.33499999999
$ws.Range("N89").Value = -86230.33499999999
$ws.Range("H122").Value = 2030
$ws.Range("J122").Value = 4600
$ws.Range("L122").Value = 13800
$ws.Range("N122").Value = -18700
$ws.Range("H134").Value = 2234.432
$ws.Range("I134").Value = 1910.425
$ws.Range("J134").Value = 5474.5
$ws.Range("K134").Value = 5731.275
$ws.Range("L134").Value = 16423.5
$ws.Range("M134").Value = -3196.275
$ws.Range("N134").Value = -21493.5
$ws.Range("H136").Value = 5980081.5
$ws.Range("I136").Value = 1283.5807
$ws.Range("K136").Value = 3850.7421
$ws.Range("M136").Value = -1300.7421

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1662.3636
$ws.Range("I92").Value = 1561.3334
$ws.Range("J92").Value = 1783.6
$ws.Range("K92").Value = 4684.0002
$ws.Range("L92").Value = 5350.799999999999
$ws.Range("M92").Value = -3436.0002
$ws.Range("N92").Value = -7846.799999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3907.6316
$ws.Range("I113").Value = 3430
$ws.Range("K113").Value = 3430
$ws.Range("M113").Value = -1260

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = 0

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2896.2712
$ws.Range("I81").Value = 2576.2307
$ws.Range("J81").Value = 3520.35
$ws.Range("K81").Value = 5152.4614
$ws.Range("L81").Value = 7040.7
$ws.Range("M81").Value = -4091.4614
$ws.Range("N81").Value = -9162.700000000001
$ws.Range("H84").Value = 2896.2712
$ws.Range("I84").Value = 2576.2307
$ws.Range("J84").Value = 3520.35
$ws.Range("K84").Value = 25762.307
$ws.Range("L84").Value = 35203.5
$ws.Range("M84").Value = -20458.307
$ws.Range("N84").Value = -45811.5
$ws.Range("H100").Value = 91820570
$ws.Range("I100").Value = 91820570
$ws.Range("K100").Value = 183641140
$ws.Range("M100").Value = -183640599
$ws.Range("H109").Value = 60000
$ws.Range("I109").Value = 60000
$ws.Range("K109").Value = 60000
$ws.Range("M109").Value = -58613
$ws.Range("H115").Value = 29999.25
$ws.Range("J115").Value = 29999.666
$ws.Range("L115").Value = 29999.666
$ws.Range("N115").Value = -33133.666
